$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells to use FV2404/FV2410 suffixes instead of old/new
$suffixMap = @{
    "Segmentname_old" = "Segmentname_FV2404"
    "Segmentgruppe_old" = "Segmentgruppe_FV2404"
    "Segment_old" = "Segment_FV2404"
    "Datenelement_old" = "Datenelement_FV2404"
    "Segment ID_old" = "Segment ID_FV2404"
    "Code_old" = "Code_FV2404"
    "Qualifier_old" = "Qualifier_FV2404"
    "Beschreibung_old" = "Beschreibung_FV2404"
    "Bedingungsausdruck_old" = "Bedingungsausdruck_FV2404"
    "Bedingung_old" = "Bedingung_FV2404"
    "Segmentname_new" = "Segmentname_FV2410"
    "Segmentgruppe_new" = "Segmentgruppe_FV2410"
    "Segment_new" = "Segment_FV2410"
    "Datenelement_new" = "Datenelement_FV2410"
    "Segment ID_new" = "Segment ID_FV2410"
    "Code_new" = "Code_FV2410"
    "Qualifier_new" = "Qualifier_FV2410"
    "Beschreibung_new" = "Beschreibung_FV2410"
    "Bedingungsausdruck_new" = "Bedingungsausdruck_FV2410"
    "Bedingung_new" = "Bedingung_FV2410"
}

for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cur = $cell.Value()
    if ($suffixMap.ContainsKey($cur)) {
        $cell.Value = $suffixMap[$cur]
    }
}

# Freeze top row
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# Create a table over the data range
$rng = $ws.Range("A1:U85")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
